$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.108.77'
$ws.Range("E2").Value = '  -5.02%  '
$ws.Range("D3").Value = '3.098.02'
$ws.Range("E3").Value = '  -6.12%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '556.39'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '160.23'
$ws.Range("E6").Value = '  -11.26%  '
$ws.Range("E7").Value = '  +0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.573'
$ws.Range("E8").Value = '  -10.47%  '
$ws.Range("D9").Value = '3.088.76'
$ws.Range("E9").Value = '  -6.57%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.67'
$ws.Range("E10").Value = '  -2.63%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.113'
$ws.Range("E11").Value = '  -10.17%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.371'
$ws.Range("E12").Value = '  -7.91%  '
$ws.Range("D13").Value = '3.639.04'
$ws.Range("E13").Value = '  -6.19%  '
$ws.Range("E14").Value = '  -1.95%  '
$ws.Range("D15").Value = '63.036.83'
$ws.Range("E15").Value = '  -5.15%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '24.18'
$ws.Range("E16").Value = '  -8.97%  '
$ws.Range("D17").Value = '3.106.73'
$ws.Range("E17").Value = '  -5.97%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000150'
$ws.Range("E18").Value = '  -8.28%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '393.11'
$ws.Range("E19").Value = '  -8.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.23'
$ws.Range("E20").Value = '  -6.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.06'
$ws.Range("E21").Value = '  -7.60%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.90'
$ws.Range("E22").Value = '  -5.34%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.997'
$ws.Range("E23").Value = '  -0.26%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.69'
$ws.Range("E24").Value = '  -0.53%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '67.16'
$ws.Range("E25").Value = '  -6.04%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.197'
$ws.Range("E26").Value = '  -4.48%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.470'
$ws.Range("E27").Value = '  -8.03%  '
$ws.Range("D28").Value = '0.0₃0985'
$ws.Range("E28").Value = '  -13.61%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.22%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.46'
$ws.Range("E30").Value = '  -9.36%  '
$ws.Range("E31").Value = '  -0.12%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.75'
$ws.Range("E32").Value = '  -8.68%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.65'
$ws.Range("E33").Value = '  -7.38%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.10'
$ws.Range("E34").Value = '  -7.07%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.69'
$ws.Range("E35").Value = '  -9.44%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '153.58'
$ws.Range("E36").Value = '  -4.01%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.08'
$ws.Range("E37").Value = '  -9.00%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.28'
$ws.Range("E38").Value = '  -10.23%  '
$ws.Range("D39").Value = '2.674.05'
$ws.Range("E39").Value = '  -6.29%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.62'
$ws.Range("E40").Value = '  -9.40%  '
$ws.Range("B41").Value = 'OKB'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '38.28'
$ws.Range("E41").Value = '  -3.63%  '
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.75'
$ws.Range("E42").Value = '  -13.47%  '
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.93'
$ws.Range("E43").Value = '  -9.55%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.687'
$ws.Range("E44").Value = '  -8.40%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0594'
$ws.Range("E45").Value = '  -7.49%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.30'
$ws.Range("E46").Value = '  -10.41%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0250'
$ws.Range("E47").Value = '  -7.35%  '
$ws.Range("B48").Value = 'FirstDigitalUSD'
$ws.Range("C48").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.00'
$ws.Range("E48").Value = '  -0.07%  '
$ws.Range("B49").Value = 'Bittensor'
$ws.Range("C49").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '277.87'
$ws.Range("E49").Value = '  -10.97%  '
$ws.Range("B50").Value = 'InjectiveProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '20.32'
$ws.Range("E50").Value = '  -10.88%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0960'
$ws.Range("E51").Value = '  -6.32%  '
